$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to remain text even when the new value looks like a
    # plain number (e.g. "233.33"), then restore the default (no explicit
    # number format / style) so the saved XML matches a plain text cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "35.042.28"
$ws.Range("E2").Value = "  +0.97%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.845.28"
$ws.Range("E3").Value = "  +2.04%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextCell "D5" "233.33"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +2.86%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - Solana
Set-TextCell "D8" "41.75"
$ws.Range("E8").Value = "  +6.07%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.22%  "

# Row 10 - Dogecoin
Set-TextCell "D10" "0.0695"
$ws.Range("E10").Value = "  +1.95%  "

# Row 11 - TRON
Set-TextCell "D11" "0.0980"
$ws.Range("E11").Value = "  -1.29%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.112.94"
$ws.Range("E12").Value = "  +2.10%  "

# Row 13 - Chainlink
Set-TextCell "D13" "11.50"
$ws.Range("E13").Value = "  +4.28%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.844.73"
$ws.Range("E14").Value = "  +2.16%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.07%  "

# Row 16 - Polkadot
Set-TextCell "D16" "4.70"
$ws.Range("E16").Value = "  +2.46%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "35.063.17"
$ws.Range("E17").Value = "  +1.14%  "

# Row 18 - Litecoin
Set-TextCell "D18" "70.00"
$ws.Range("E18").Value = "  +0.72%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0791"
$ws.Range("E19").Value = "  +0.59%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "240.57"
$ws.Range("E20").Value = "  +0.24%  "

# Row 21 - Avalanche
Set-TextCell "D21" "12.17"
$ws.Range("E21").Value = "  +1.88%  "

# Row 22 - Uniswap
Set-TextCell "D22" "4.78"
$ws.Range("E22").Value = "  +2.58%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.01%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +3.17%  "

# Row 25 - Monero
Set-TextCell "D25" "172.11"
$ws.Range("E25").Value = "  +0.02%  "

# Row 26 - Cosmos
Set-TextCell "D26" "7.91"

# Row 27 - EthereumClassic
Set-TextCell "D27" "17.52"
$ws.Range("E27").Value = "  +1.85%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +3.61%  "

# Row 29 - PancakeSwap
Set-TextCell "D29" "1.73"
$ws.Range("E29").Value = "  +11.17%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.06%  "

# Row 31 - Hedera
Set-TextCell "D31" "0.0556"
$ws.Range("E31").Value = "  +1.95%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -1.25%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextCell "D33" "3.95"
$ws.Range("E33").Value = "  -0.75%  "

# Row 34 - WEMIXToken
$ws.Range("E34").Value = "  +23.66%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +10.73%  "

# Row 36 - ImmutableX
Set-TextCell "D36" "0.765"
$ws.Range("E36").Value = "  +9.25%  "

# Row 37 - TrustWalletToken
$ws.Range("E37").Value = "  -2.53%  "

# Row 38 - ARBITRUM
Set-TextCell "D38" "1.07"
$ws.Range("E38").Value = "  +11.35%  "

# Row 39 - was VeChain, now Aave
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D39" "90.25"
$ws.Range("E39").Value = "  -1.43%  "

# Row 40 - was Aave, now VeChain
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D40" "0.0201"
$ws.Range("E40").Value = "  +4.62%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.344.85"
$ws.Range("E41").Value = "  +1.80%  "

# Row 42 - InjectiveProtocol
Set-TextCell "D42" "14.63"
$ws.Range("E42").Value = "  +2.61%  "

# Row 43 - RenderToken
Set-TextCell "D43" "2.29"
$ws.Range("E43").Value = "  +3.63%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  +3.91%  "

# Row 45 - HuobiToken
$ws.Range("E45").Value = "  -4.08%  "

# Row 46 - was Gas, now Kaspa
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D46" "0.0532"
$ws.Range("E46").Value = "  +3.80%  "

# Row 47 - was Kaspa, now FraxShare
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D47" "6.34"
$ws.Range("E47").Value = "  +1.38%  "

# Row 48 - was FraxShare, now Gas
$ws.Range("B48").Value = "Gas"
$ws.Range("C48").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
Set-TextCell "D48" "11.72"
$ws.Range("E48").Value = "  +74.34%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.027.99"
$ws.Range("E49").Value = "  +1.54%  "

# Row 50 - THORChain
Set-TextCell "D50" "3.43"
$ws.Range("E50").Value = "  +16.02%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  +0.57%  "
